# Rarres2 -> Ccrl2 LR-pair sheet: refresh with new TPM-based NATMI stats.
# Sending/target cluster set grows from {ECs, FAPs, MuSCs} to {ECs, FAPs, MuSCs, Resolving-Mac}
# (Resolving-Mac moves from a target-only cluster to also being a sending cluster).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")
$ligand   = "Rarres2"
$receptor = "Ccrl2"

# Ligand-side stats (detection/expression), keyed by sending cluster: E,F,G,H,I,J
$ef = @{}
$ef["ECs"] = @(3, 1)
$ef["FAPs"] = @(3, 1)
$ef["MuSCs"] = @(3, 1)
$ef["Resolving-Mac"] = @(1, 0.3333333333333333)
$ligandStats = @{}
$ligandStats["ECs"] = @(1.458525, 4.375575, 0.02525273220298681, 0.02525273220298681)
$ligandStats["FAPs"] = @(14.58505066666667, 43.755152, 0.84033742450786, 0.8403374245078601)
$ligandStats["MuSCs"] = @(7.712822, 23.138466, 0.1335389029981009, 0.133538902998101)
$ligandStats["Resolving-Mac"] = @(0.05030299999999999, 0.150909, 0.0008709402910521559, 0.000870940291052156)

# Receptor-side stats (detection/expression), keyed by target cluster: K,L,M,N,O,P
$kl = @{}
$kl["ECs"] = @(3, 1)
$kl["FAPs"] = @(2, 0.6666666666666666)
$kl["MuSCs"] = @(3, 1)
$kl["Resolving-Mac"] = @(3, 1)
$receptorStats = @{}
$receptorStats["ECs"] = @(4.306835666666667, 12.920507, 0.0474068797554582, 0.04740687975545821)
$receptorStats["FAPs"] = @(0.464105, 1.392315, 0.005108569639466994, 0.005108569639466995)
$receptorStats["MuSCs"] = @(1.661774, 4.985322, 0.01829174045540476, 0.01829174045540476)
$receptorStats["Resolving-Mac"] = @(84.41561133333333, 253.246834, 0.92919281014967, 0.9291928101496701)

# Edge weight / derived-specificity stats, keyed by "sending|target": Q,R,S,T
$edgeStats = @{}
$edgeStats["ECs|ECs"] = @(6.281627490725, 56.534647416525, 0.001197153239043783, 0.001197153239043783)
$edgeStats["ECs|FAPs"] = @(0.6769087451249999, 6.092178706124999, 0.0001290053410457689, 0.0001290053410457689)
$edgeStats["ECs|MuSCs"] = @(2.42373892335, 21.81365031015, 0.0004619164232468764, 0.0004619164232468766)
$edgeStats["ECs|Resolving-Mac"] = @(123.12227951995, 1108.10051567955, 0.02346465719965038, 0.02346465719965039)
$edgeStats["FAPs|ECs"] = @(209.0342789383112, 1881.3085104448, 0.03983777523765555, 0.03983777523765556)
$edgeStats["FAPs|FAPs"] = @(22.52555275733334, 202.729974816, 0.004292922253748741, 0.004292922253748742)
$edgeStats["FAPs|MuSCs"] = @(80.65497658453334, 725.8947892608, 0.01537123406406107, 0.01537123406406107)
$edgeStats["FAPs|Resolving-Mac"] = @(4097.151090015289, 36874.3598101376, 0.7808354929523946, 0.7808354929523948)
$edgeStats["MuSCs|ECs"] = @(33.21785688025134, 298.960711922262, 0.006330662717106767, 0.00633066271710677)
$edgeStats["MuSCs|FAPs"] = @(3.57955925431, 32.21603328879, 0.0006821927855438264, 0.0006821927855438266)
$edgeStats["MuSCs|MuSCs"] = @(12.816967066228, 115.352703596052, 0.002442658954340735, 0.002442658954340736)
$edgeStats["MuSCs|Resolving-Mac"] = @(651.0825842351827, 5859.743258116644, 0.1240833885411096, 0.1240833885411096)
$edgeStats["Resolving-Mac|ECs"] = @(0.2166467545403333, 1.949820790863, 0.00004128856165209332, 0.00004128856165209333)
$edgeStats["Resolving-Mac|FAPs"] = @(0.023345873815, 0.210112864335, 0.000004449259128657591, 0.000004449259128657592)
$edgeStats["Resolving-Mac|MuSCs"] = @(0.083592217522, 0.7523299576979999, 0.00001593101375608071, 0.00001593101375608072)
$edgeStats["Resolving-Mac|Resolving-Mac"] = @(4.246358496900666, 38.21722647210599, 0.0008092714565153242, 0.0008092714565153244)

$r = 2
foreach ($send in $clusters) {
    foreach ($target in $clusters) {
        $ef_vals = $ef[$send]
        $lig = $ligandStats[$send]
        $kl_vals = $kl[$target]
        $rec = $receptorStats[$target]
        $edge = $edgeStats["$send|$target"]

        $ws.Cells.Item($r, 1).Value = $send
        $ws.Cells.Item($r, 2).Value = $ligand
        $ws.Cells.Item($r, 3).Value = $receptor
        $ws.Cells.Item($r, 4).Value = $target
        $ws.Cells.Item($r, 5).Value = $ef_vals[0]
        $ws.Cells.Item($r, 6).Value = $ef_vals[1]
        $ws.Cells.Item($r, 7).Value = $lig[0]
        $ws.Cells.Item($r, 8).Value = $lig[1]
        $ws.Cells.Item($r, 9).Value = $lig[2]
        $ws.Cells.Item($r, 10).Value = $lig[3]
        $ws.Cells.Item($r, 11).Value = $kl_vals[0]
        $ws.Cells.Item($r, 12).Value = $kl_vals[1]
        $ws.Cells.Item($r, 13).Value = $rec[0]
        $ws.Cells.Item($r, 14).Value = $rec[1]
        $ws.Cells.Item($r, 15).Value = $rec[2]
        $ws.Cells.Item($r, 16).Value = $rec[3]
        $ws.Cells.Item($r, 17).Value = $edge[0]
        $ws.Cells.Item($r, 18).Value = $edge[1]
        $ws.Cells.Item($r, 19).Value = $edge[2]
        $ws.Cells.Item($r, 20).Value = $edge[3]

        $r = $r + 1
    }
}
